# "Revised Answers to questions."
#
# The source document repeats the placeholder run "[]" (and a few
# "[Spec. Doc. ...]" citations) across many answer paragraphs, so every
# edit below is scoped to one specific paragraph (located by its stable
# index in $d.Paragraphs) and the paragraph's own text is asserted
# against an anchor substring first -- this guarantees Find/Execute's
# wdReplaceAll (4th-to-last arg = 2) only ever touches the single
# occurrence intended, never a look-alike elsewhere in the document.

$d = $word.ActiveDocument

function Edit-Paragraph($paraIndex, $anchorSubstring, $findText, $replaceText) {
    $p = $d.Paragraphs($paraIndex)
    $t = $p.Range.Text
    if (-not $t.Contains($anchorSubstring)) {
        throw "Paragraph $paraIndex does not contain expected anchor '$anchorSubstring' (actual: '$t')"
    }
    $r = $p.Range
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace of '$findText' failed in paragraph $paraIndex (text: '$t')"
    }
}

# 1. "The file name from the second file, ... []" -> "[Spec. Doc. 4.3]"  (bold, blue)
Edit-Paragraph 6 "manually entered data" "[]" "[Spec. Doc. 4.3]"

# 2. "Should look as follows: fall2015, spring2014,... []" -> "[Spec. doc. 4.3]"
Edit-Paragraph 10 "Should look as follows" "[]" "[Spec. doc. 4.3]"

# 3. "Would like the option to print the file  [spec. Doc. 2.2/3.2]" -> 2.2 becomes 2.3
Edit-Paragraph 14 "Would like the option to print the file" "2.2" "2.3"

# 4. "A description  []" -> "[Spec. Doc. 4.5]"
Edit-Paragraph 27 "A description" "[]" "[Spec. Doc. 4.5]"

# 5. "Yes, by type []" -> "[Spec. Doc. 4.5]"
Edit-Paragraph 29 "Yes, by type" "[]" "[Spec. Doc. 4.5]"

# 6. "All numbers, no letters  []" -> "[Spec. Doc.4.3]" (no space before 4.3)
Edit-Paragraph 31 "All numbers, no letters" "[]" "[Spec. Doc.4.3]"

# 7. "Specify which portion was incorrect. []" -> "[Spec. Doc. 4.1]"
Edit-Paragraph 35 "Specify which portion was incorrect" "[]" "[Spec. Doc. 4.1]"

# 8. "Yes  [Spec. Doc. 4.1]" -> "[Spec. Doc. 2.3/3.2/4.1]"
Edit-Paragraph 39 "[Spec. Doc. 4.1]" "[Spec. Doc. 4.1]" "[Spec. Doc. 2.3/3.2/4.1]"

# 9. "Yes [Spec doc, 4.1]" -> "[Spec doc, 2.3/3.2/4.1]"
Edit-Paragraph 41 "[Spec doc, 4.1]" "[Spec doc, 4.1]" "[Spec doc, 2.3/3.2/4.1]"

# 10. "Yes [Spec. Doc. 2.2/3.2/4.1]" -> 2.2 becomes 2.3
Edit-Paragraph 43 "[Spec. Doc. 2.2/3.2/4.1]" "2.2" "2.3"

# 11. "Yes []" (re: viewing previously generated exam schedules) -> "[Spec. Doc. 3.2/2.3]"
Edit-Paragraph 47 "Yes" "[]" "[Spec. Doc. 3.2/2.3]"

# 12. "Do you want the university's logo in the program?" -> append " (Add Screen shots)"
$pLogo = $d.Paragraphs(60)
if (-not $pLogo.Range.Text.Contains("university")) {
    throw "Paragraph 60 does not contain expected anchor 'university' (actual: '$($pLogo.Range.Text)')"
}
$pLogo.Range.InsertAfter(" (Add Screen shots)") | Out-Null

# 13. Last paragraph: "AMBIGUOUS Optional views [Spec. Doc. 3.1.1]" -> "[Spec. Doc. 3.2.1.1.1]"
Edit-Paragraph 67 "AMBIGUOUS Optional views" "[Spec. Doc. 3.1.1]" "[Spec. Doc. 3.2.1.1.1]"
